# Update simulation results for Case_1_4 / res_line / pl_mw (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(1.219803940679071, 1.069997499266549, 0.9775436500530645, 0.9397517899408854, 0.9334695415156489, 0.9770344427063833, 1.168250105383493, 1.539389405566169, 1.809644099351658, 1.932050532542917, 1.978324344473322, 1.968361989860625, 1.935859098657602, 1.915939811075305, 1.801633628443199, 1.731372246120657, 1.690909632147168, 1.677201143001923, 1.738856893516015, 1.945408141631617, 2.079940855282643, 2.008181005422728, 1.735473297350438, 1.439406031741498)
for ($i = 0; $i -lt $BValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $BValues[$i]
}

$DValues = @(0.00328276805978156, 0.003117985810526136, 0.003018255499195988, 0.002977983059874845, 0.002971318247256605, 0.003017710873152168, 0.00322565360309035, 0.003644731220816766, 0.003959326331720803, 0.004103863068408486, 0.004158796915575635, 0.004146957056267553, 0.004108378505360832, 0.00408477407517438, 0.003949908852845851, 0.003867535620365459, 0.003820291157840217, 0.003804318212985436, 0.003876290507584201, 0.004119704540567426, 0.00427995902020939, 0.004194322508695336, 0.003872332070251616, 0.003530170004353295)
for ($i = 0; $i -lt $DValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $DValues[$i]
}

$EValues = @(0.4354855266417559, 0.3796236440278449, 0.3454298728940586, 0.3315199486717546, 0.3292116194363217, 0.3452421832880503, 0.4162010468717057, 0.5562976445450403, 0.6599630079396519, 0.7073182102084701, 0.7252814645633805, 0.7214113434755376, 0.7087954249356869, 0.7010719049055893, 0.6568724169760856, 0.629809763319642, 0.6142624613427046, 0.6090015003013463, 0.6326887045236589, 0.7125001682736638, 0.7648426698322197, 0.7368890891373212, 0.6313870999755409, 0.518280098549738)
for ($i = 0; $i -lt $EValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $EValues[$i]
}

$FValues = @(0.4268455756567278, 0.3933003553665202, 0.3729420886923549, 0.3647054950681934, 0.3633413997772692, 0.3728307665163442, 0.4152293778963099, 0.5002909339772401, 0.5639956345128354, 0.5932485116357213, 0.604365721277091, 0.6019696543234119, 0.5941623311389179, 0.5893853151532795, 0.5620894304974513, 0.5454146261146775, 0.5358494138644119, 0.5326151957723084, 0.5471870247896504, 0.5964544477321994, 0.6288856518111032, 0.6115551107116914, 0.5463856568671588, 0.4770700030462791)
for ($i = 0; $i -lt $FValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $FValues[$i]
}

$GValues = @(0.3367752227119354, 0.3012008603001277, 0.2794993337179079, 0.2706907303260948, 0.269230162170615, 0.2793803970864275, 0.3244794999739042, 0.4140703990656505, 0.4806440831150098, 0.5111043360739984, 0.5226648403903766, 0.5201739184791165, 0.5120549037030742, 0.5070851600557091, 0.4786570196494324, 0.4612626808594769, 0.4512744032264209, 0.4478953566152768, 0.4631126271809762, 0.5144389481118878, 0.5481347357242612, 0.5301366292119667, 0.4622762292457878, 0.3897053291205168)
for ($i = 0; $i -lt $GValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $GValues[$i]
}

$HValues = @(0.3434920082859492, 0.3306581927846537, 0.3229880749787242, 0.3199150623406126, 0.3194079646665671, 0.3229464183999653, 0.3390232442816625, 0.3722245405283218, 0.3976543754028228, 0.4094518639923876, 0.4139524800987431, 0.4129817139520924, 0.409821466670877, 0.4078900483665961, 0.3968880122874339, 0.3901974874650023, 0.3863708487733675, 0.3850789144428575, 0.3909074710753089, 0.4107488063632729, 0.4239096552953754, 0.4168677050664655, 0.3905864257188227, 0.3630617902832967)
for ($i = 0; $i -lt $HValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 8).Value = $HValues[$i]
}

$IValues = @(3.961193417866525, 3.481483530427994, 3.186307868167489, 3.065876933603903, 3.04587117873993, 3.184684259858585, 3.795928923952232, 4.989014804984947, 5.861568061238756, 6.25753114434093, 6.407322560669854, 6.375069214559744, 6.269857673435979, 6.205392518134659, 5.835670441680065, 5.608601612490133, 5.477907426105048, 5.433641513037287, 5.632782903638144, 6.30076505272973, 6.736443306619492, 6.503998686294949, 5.621851005971507, 4.666914883151151)
for ($i = 0; $i -lt $IValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 9).Value = $IValues[$i]
}

$OValues = @(1.288635611322178, 1.186000399322808, 1.123707417123967, 1.098503372127993, 1.094329137056548, 1.12336677637748, 1.253095441307636, 1.513334097579957, 1.70823048577472, 1.797728296575315, 1.831741466706262, 1.824410655378244, 1.800524119720649, 1.785908904299276, 1.702398612130764, 1.651383722660285, 1.622120123892898, 1.612225456820568, 1.656806175146926, 1.807536841490844, 1.906761666784462, 1.853737639943063, 1.654354480355209, 1.442293097308436)
for ($i = 0; $i -lt $OValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 15).Value = $OValues[$i]
}
